# Applies the LOQ4064.docx content update:
#  - Créditos-trabalho: 4 -> 0
#  - Carga horária: 150 h -> 30 h
#  - Ativação: 01/01/2016 -> 01/01/2025
#  - Collapse multi-run/line-break paragraphs (Objetivos EN, Programa PT, Programa EN)
#    into single runs with the manual line breaks removed (text simply concatenated).

$d = $word.ActiveDocument
$vtab = [char]11

# --- Simple header field updates -------------------------------------------------

$d.Content.Find.Execute("Créditos-trabalho: 4", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Créditos-trabalho: 0", 2) | Out-Null

$d.Content.Find.Execute("Carga horária: 150 h", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Carga horária: 30 h", 2) | Out-Null

$d.Content.Find.Execute("Ativação: 01/01/2016", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2025", 2) | Out-Null

# --- Objetivos (English, italic run) : merge the two <w:t> runs split by <w:br/> ----

$objEn1 = "1 - Consolidation and application of knowledge acquired in each of the specific areas of the Chemical Engineering degree. "
$objEn2 = "2 - Integration of knowledge of Chemical Engineering"

$searchObjEn = $objEn1 + $vtab + $objEn2
$replaceObjEn = $objEn1 + $objEn2

$d.Content.Find.Execute($searchObjEn, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replaceObjEn, 2) | Out-Null

# --- Programa (Portuguese) : merge the four <w:t> runs split by <w:br/> ------------

$pt1 = "1 - Diagramas para estudos de processos químicos: diagramas de bloco; Fluxogramas de processo (PFD); Fluxogramas de instrumentação e tubulação (P&ID)."
$pt2 = "2  Estrutura e síntese de processos químicos industriais: Hierarquia no planejamento de processos; Etapa 1- Descontínuo ou contínuo; Etapa 2 - Estrutura de entrada/saída de processo; Etapa 3- Estrutura de reciclo; "
$pt3 = "3  Análise de desempenho de processos químicos: Modelo de entrada e saída; Ferramentas para a avaliação de processos."
$pt4 = "4  Estudo de planta química industrial."

$searchPrograma = $pt1 + $vtab + $pt2 + $vtab + $pt3 + $vtab + $pt4
$replacePrograma = $pt1 + $pt2 + $pt3 + $pt4

$d.Content.Find.Execute($searchPrograma, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replacePrograma, 2) | Out-Null

# --- Programa (English, italic run) : merge the four <w:t> runs split by <w:br/> ---

$en1 = "1 - Diagrams for Understanding Chemical Processes: Block Flow Diagrams; Process Flow Diagram (PFD); Piping and Instrumentation Diagram (P&ID)."
$en2 = "2 - The Structure and Synthesis of Process Flow Diagrams:  Hierarchy of Process Design; Step 1 - Batch versus Continuous Process; Step 2 - The Input/Output Structure of the Process; Step 3 - The Recycle Structure of the Process"
$en3 = "3 - Analysis of process performance: Process Input/Output Models; Tools for evaluating process performance."
$en4 = "4 - Industrial chemical plant study."

$searchProgramaEn = $en1 + $vtab + $en2 + $vtab + $en3 + $vtab + $en4
$replaceProgramaEn = $en1 + $en2 + $en3 + $en4

$d.Content.Find.Execute($searchProgramaEn, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replaceProgramaEn, 2) | Out-Null

Write-Host "done"
